# Applies the "updated Medical Decision problem" edit to MDDMparameters.xlsx
#
# Summary of content changes:
#  - Instructions!A2: the explanatory example text now references question 4.16 /
#    4.17 a) / 4.17 b) instead of question 2 / 3.a / 3.b (bold formatting on
#    "parameters2", "truth_type" and the "Possible values..." phrase preserved)
#  - parameters1!C3 (Sens sigma_0): 0.19 -> 0.09
#  - parameters1!C6 (PA sigma_0):   0.21 -> 0.11
#  - parameters2!B2 (sigma_w):      0.5  -> 5
#  - parameters2!B5 (theta_start):  0    -> 1
#  - parameters2!B6 (theta_end):    2.1  -> 1
#  - parameters2!B8 (truth_type):   "fixed_uniform" -> "known"

$wb = $excel.ActiveWorkbook

$wsInstructions = $wb.Worksheets.Item("Instructions")
$wsParameters1  = $wb.Worksheets.Item("parameters1")
$wsParameters2  = $wb.Worksheets.Item("parameters2")

# ---------------------------------------------------------------------------
# 1) Instructions sheet: update the worked-example question numbers in A2,
#    while preserving the bold runs already present in that rich-text cell.
#    Doing this BEFORE the parameters2!B8 edit below keeps the new "known"
#    shared string appended last, matching how the workbook was re-saved.
# ---------------------------------------------------------------------------
$a2 = $wsInstructions.Range("A2")
$originalText = $a2.Value2

$updatedText = $originalText.Replace('question 2,', 'question 4.16,').Replace('question 3.a,', 'question 4.17 a),').Replace('question 3.b,', 'question 4.17 b),')

$a2.Value = $updatedText

function Set-BoldRun($Cell, $SearchText, $SearchStart) {
    $current = $Cell.Value2
    $zeroBasedStart = $SearchStart - 1
    $foundAt = $current.IndexOf($SearchText, $zeroBasedStart)
    if ($foundAt -lt 0) {
        return $SearchStart
    }
    $oneBasedAt = $foundAt + 1
    $Cell.Characters($oneBasedAt, $SearchText.Length).Font.Bold = $true
    return $oneBasedAt + $SearchText.Length
}

$cursor = 1
$cursor = Set-BoldRun $a2 "parameters2" $cursor
$cursor = Set-BoldRun $a2 "truth_type" $cursor
$cursor = Set-BoldRun $a2 'Possible values for "truth_type" are "known", "fixed_uniform", "prior_uniform" or "normal"' $cursor

# ---------------------------------------------------------------------------
# 2) parameters1 sheet: adjust the sigma_0 values for Sens and PA
# ---------------------------------------------------------------------------
$wsParameters1.Range("C3").Value = 0.09
$wsParameters1.Range("C6").Value = 0.11

# ---------------------------------------------------------------------------
# 3) parameters2 sheet: updated algorithm parameters
# ---------------------------------------------------------------------------
$wsParameters2.Range("B2").Value = 5
$wsParameters2.Range("B5").Value = 1
$wsParameters2.Range("B6").Value = 1
$wsParameters2.Range("B8").Value = "known"
